$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "This is a contract made between the Seller, {{ seller.name }}, social insurance number {{ seller.sin }}" "{{ sellers[0].name.full() }}, social insurance number {{ sellers[0].social_insurance_number }}"
Replace-Text "and the Buyer, {{ buyer.name }}" "{{ buyers[0].name.full() }}"
Replace-Text "social insurance number {{ buyer.sin }}," "social insurance number {{ buyers[0].social_insurance_number }}, "
Replace-Text "for the sale of Seller’s vehicle:" "for the sale of Seller’s "
Replace-Text "Year: {{ vehicle.year }}" "{{ vehicle.year }}"
Replace-Text "Make: {{ vehicle.make }}" ": {{ vehicle.make }}"
Replace-Text "Model: {{ vehicle.model }}" ": {{ vehicle.model }}"
Replace-Text "Exterior color: {{ vehicle.exterior_color }}" ": {{ vehicle.exterior_color }}"
Replace-Text "Interior color: {{ vehicle.interior_color }}" ": {{ vehicle.interior_color }} "
Replace-Text "The VIN number is {{ vehicle.vin }}" "{{ vehicle.vin_number }}"
Replace-Text "and the odometer reads {{ vehicle.odometer }}" "{{ vehicle.odometer_reading }}"
Replace-Text "miles as of {{ sale_date }}." "{{ signature_date }}"
Replace-Text "The date of sale is {{ sale_date }}." "{{ signature_date }}"
Replace-Text "Buyer agrees to pay to Seller the purchase price of {{ sale_price }}" "{{ sale.purchase_price }}"
Replace-Text "to be paid in {{ payment_type }}." "{{ sale.payment_method }}"
Replace-Text "Seller will provide the Buyer with the vehicle’s title and {{ required_documents }}." "{{ vehicle.documents_provided }}"
Replace-Text "DATE: {{ sale_date }}" "DATE: `${today}"
